# extraction and standardization code for healthcare
#
# Applies:
#  1. Updates the "amenity+building" notes text (row 39 / F39, also reused by
#     F40:F43) to reflect that the columns have already been dropped.
#  2. Adds a new note for the Healthcare section describing the
#     'amenity' + 'building' + 'healthcare' combination, and stamps the
#     Healthcare rows (27-38) with the same "needs review" highlight /
#     dictionary-key / note pattern used by the other sections.
#  3. Adds the same dictionary-key / note pair to F40:F43 (kindergarten,
#     library, school, university) which previously only had the
#     dictionary key (E) filled in.
#  4. Resets the sheet selection (no more scrolled-down view pinned on
#     G41; the workbook now opens selected at I3 at the top).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Infrastructure categorization")

$oldNote = "Contains column 'amenity' and 'building' which are combined into one column. Drop original columns. Decide on reclassification."
$updatedNote = "Contains column 'amenity' and 'building' which are combined into one column. Dropped original columns. Still need to decide on reclassification."
$healthcareNote = "Contains column 'amenity', 'building' and 'healthcare' which are combined into one column. Dropped original columns. Still need to decide on reclassification."

# 1. Update the existing note text in place (keeps the same shared string,
#    referenced by F39, and now also by F40:F43).
$ws.Range("F39").Value = $updatedNote

# Propagate the (updated) note + its dictionary key to the remaining
# Education rows that were missing them.
$eduRows = 40,41,42,43
foreach ($r in $eduRows) {
    $ws.Cells.Item($r, 6).Value = $updatedNote
}

# 2. Healthcare section (rows 27-38): add the OSM-flex "needs review" fill
#    in column D (copy formatting from D17, which already carries it),
#    plus the dictionary key (E) and the new healthcare note (F).
$fillSource = $ws.Range("D17")
$fillSource.Copy()

$healthRows = 27,28,29,30,31,32,33,34,35,36,37,38
foreach ($r in $healthRows) {
    $dst = $ws.Cells.Item($r, 4)
    $dst.PasteSpecial(-4122)
    $ws.Cells.Item($r, 5).Value = "healthcare"
    $ws.Cells.Item($r, 6).Value = $healthcareNote
}

$excel.CutCopyMode = $false

# 3. Reset view: select I3, clear the old scrolled/selected state (G41).
$ws.Activate()
$ws.Range("I3").Select()
